$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.8238633333333333
$ws.Range("H2").Value = 2.47159
$ws.Range("I2").Value = 0.2181573957783601
$ws.Range("J2").Value = 0.2181573957783601
$ws.Range("M2").Value = 1.021465
$ws.Range("N2").Value = 3.064395
$ws.Range("O2").Value = 0.03918894356403102
$ws.Range("P2").Value = 0.03918894356403101
$ws.Range("Q2").Value = 0.8415475597833333
$ws.Range("R2").Value = 7.57392803805
$ws.Range("S2").Value = 0.008549357871234134
$ws.Range("T2").Value = 0.008549357871234132
$ws.Range("G3").Value = 0.8238633333333333
$ws.Range("H3").Value = 2.47159
$ws.Range("I3").Value = 0.2181573957783601
$ws.Range("J3").Value = 0.2181573957783601
$ws.Range("O3").Value = 0.3981426681675393
$ws.Range("P3").Value = 0.3981426681675393
$ws.Range("Q3").Value = 8.549758181017777
$ws.Range("R3").Value = 76.94782362916
$ws.Range("S3").Value = 0.08685776763567819
$ws.Range("T3").Value = 0.08685776763567818
$ws.Range("G4").Value = 0.8238633333333333
$ws.Range("H4").Value = 2.47159
$ws.Range("I4").Value = 0.2181573957783601
$ws.Range("J4").Value = 0.2181573957783601
$ws.Range("M4").Value = 4.368617666666666
$ws.Range("N4").Value = 13.105853
$ws.Range("O4").Value = 0.1676038936153748
$ws.Range("P4").Value = 0.1676038936153748
$ws.Range("Q4").Value = 3.599143912918888
$ws.Range("R4").Value = 32.39229521627
$ws.Range("S4").Value = 0.0365640289534435
$ws.Range("T4").Value = 0.0365640289534435
$ws.Range("G5").Value = 0.8238633333333333
$ws.Range("H5").Value = 2.47159
$ws.Range("I5").Value = 0.2181573957783601
$ws.Range("J5").Value = 0.2181573957783601
$ws.Range("M5").Value = 10.29740833333333
$ws.Range("N5").Value = 30.892225
$ws.Range("O5").Value = 0.3950644946530549
$ws.Range("P5").Value = 0.3950644946530548
$ws.Range("Q5").Value = 8.483657154194443
$ws.Range("R5").Value = 76.35291438775
$ws.Range("S5").Value = 0.08618624131800433
$ws.Range("T5").Value = 0.08618624131800431
$ws.Range("I6").Value = 0.3205322899584435
$ws.Range("J6").Value = 0.3205322899584435
$ws.Range("M6").Value = 1.021465
$ws.Range("N6").Value = 3.064395
$ws.Range("O6").Value = 0.03918894356403102
$ws.Range("P6").Value = 0.03918894356403101
$ws.Range("Q6").Value = 1.236461250758333
$ws.Range("R6").Value = 11.128151256825
$ws.Range("S6").Value = 0.01256132182163107
$ws.Range("T6").Value = 0.01256132182163107
$ws.Range("I7").Value = 0.3205322899584435
$ws.Range("J7").Value = 0.3205322899584435
$ws.Range("O7").Value = 0.3981426681675393
$ws.Range("P7").Value = 0.3981426681675393
$ws.Range("Q7").Value = 12.56190998510444
$ws.Range("S7").Value = 0.1276175811579061
$ws.Range("T7").Value = 0.127617581157906
$ws.Range("I8").Value = 0.3205322899584435
$ws.Range("J8").Value = 0.3205322899584435
$ws.Range("M8").Value = 4.368617666666666
$ws.Range("N8").Value = 13.105853
$ws.Range("O8").Value = 0.1676038936153748
$ws.Range("P8").Value = 0.1676038936153748
$ws.Range("Q8").Value = 5.288117032117222
$ws.Range("R8").Value = 47.593053289055
$ws.Range("S8").Value = 0.05372245982648744
$ws.Range("T8").Value = 0.05372245982648744
$ws.Range("I9").Value = 0.3205322899584435
$ws.Range("J9").Value = 0.3205322899584435
$ws.Range("M9").Value = 10.29740833333333
$ws.Range("N9").Value = 30.892225
$ws.Range("O9").Value = 0.3950644946530549
$ws.Range("P9").Value = 0.3950644946530548
$ws.Range("Q9").Value = 12.46478967698611
$ws.Range("R9").Value = 112.183107092875
$ws.Range("S9").Value = 0.1266309271524189
$ws.Range("T9").Value = 0.1266309271524189
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.316433
$ws.Range("H10").Value = 0.9492990000000001
$ws.Range("I10").Value = 0.0837908381466997
$ws.Range("J10").Value = 0.0837908381466997
$ws.Range("M10").Value = 1.021465
$ws.Range("N10").Value = 3.064395
$ws.Range("O10").Value = 0.03918894356403102
$ws.Range("P10").Value = 0.03918894356403101
$ws.Range("Q10").Value = 0.3232252343450001
$ws.Range("R10").Value = 2.909027109105001
$ws.Range("S10").Value = 0.003283674427313872
$ws.Range("T10").Value = 0.003283674427313872
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.316433
$ws.Range("H11").Value = 0.9492990000000001
$ws.Range("I11").Value = 0.0837908381466997
$ws.Range("J11").Value = 0.0837908381466997
$ws.Range("O11").Value = 0.3981426681675393
$ws.Range("P11").Value = 0.3981426681675393
$ws.Range("Q11").Value = 3.283828180030667
$ws.Range("R11").Value = 29.554453620276
$ws.Range("S11").Value = 0.03336070786772145
$ws.Range("T11").Value = 0.03336070786772145
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.316433
$ws.Range("H12").Value = 0.9492990000000001
$ws.Range("I12").Value = 0.0837908381466997
$ws.Range("J12").Value = 0.0837908381466997
$ws.Range("M12").Value = 4.368617666666666
$ws.Range("N12").Value = 13.105853
$ws.Range("O12").Value = 0.1676038936153748
$ws.Range("P12").Value = 0.1676038936153748
$ws.Range("Q12").Value = 1.382374794116333
$ws.Range("R12").Value = 12.441373147047
$ws.Range("S12").Value = 0.01404367072268255
$ws.Range("T12").Value = 0.01404367072268255
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.316433
$ws.Range("H13").Value = 0.9492990000000001
$ws.Range("I13").Value = 0.0837908381466997
$ws.Range("J13").Value = 0.0837908381466997
$ws.Range("M13").Value = 10.29740833333333
$ws.Range("N13").Value = 30.892225
$ws.Range("O13").Value = 0.3950644946530549
$ws.Range("P13").Value = 0.3950644946530548
$ws.Range("Q13").Value = 3.258439811141667
$ws.Range("R13").Value = 29.325958300275
$ws.Range("S13").Value = 0.03310278512898183
$ws.Range("T13").Value = 0.03310278512898183
$ws.Range("G14").Value = 1.425688333333333
$ws.Range("H14").Value = 4.277065
$ws.Range("I14").Value = 0.3775194761164967
$ws.Range("J14").Value = 0.3775194761164967
$ws.Range("M14").Value = 1.021465
$ws.Range("N14").Value = 3.064395
$ws.Range("O14").Value = 0.03918894356403102
$ws.Range("P14").Value = 0.03918894356403101
$ws.Range("Q14").Value = 1.456290733408333
$ws.Range("R14").Value = 13.106616600675
$ws.Range("S14").Value = 0.01479458944385195
$ws.Range("T14").Value = 0.01479458944385194
$ws.Range("G15").Value = 1.425688333333333
$ws.Range("H15").Value = 4.277065
$ws.Range("I15").Value = 0.3775194761164967
$ws.Range("J15").Value = 0.3775194761164967
$ws.Range("O15").Value = 0.3981426681675393
$ws.Range("P15").Value = 0.3981426681675393
$ws.Range("Q15").Value = 14.79528217645111
$ws.Range("R15").Value = 133.15753958806
$ws.Range("S15").Value = 0.1503066115062336
$ws.Range("T15").Value = 0.1503066115062336
$ws.Range("G16").Value = 1.425688333333333
$ws.Range("H16").Value = 4.277065
$ws.Range("I16").Value = 0.3775194761164967
$ws.Range("J16").Value = 0.3775194761164967
$ws.Range("M16").Value = 4.368617666666666
$ws.Range("N16").Value = 13.105853
$ws.Range("O16").Value = 0.1676038936153748
$ws.Range("P16").Value = 0.1676038936153748
$ws.Range("Q16").Value = 6.228287240160555
$ws.Range("R16").Value = 56.054585161445
$ws.Range("S16").Value = 0.06327373411276135
$ws.Range("T16").Value = 0.06327373411276135
$ws.Range("G17").Value = 1.425688333333333
$ws.Range("H17").Value = 4.277065
$ws.Range("I17").Value = 0.3775194761164967
$ws.Range("J17").Value = 0.3775194761164967
$ws.Range("M17").Value = 10.29740833333333
$ws.Range("N17").Value = 30.892225
$ws.Range("O17").Value = 0.3950644946530549
$ws.Range("P17").Value = 0.3950644946530548
$ws.Range("Q17").Value = 14.68089492440278
$ws.Range("R17").Value = 132.128054319625
$ws.Range("S17").Value = 0.1491445410536498
$ws.Range("T17").Value = 0.1491445410536497
